$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated on" timestamp (pav:lastUpdatedOn value in B23)
$ws.Range("B23").Value = "2022-06-28T06:38:59+00:00"

# The regenerated vocabulary dropped the stray duplicate row 73
# ("concentration of soluble factors in human blood samples", which was
# already fully defined further down). Deleting it shifts every following
# row up by one -- Excel keeps the dimension / row indices in sync
# automatically.
$ws.Rows("73").Delete()

# The sequential "incentive:NNNN" identifiers in column A are regenerated
# from scratch against the new row order, so every remaining auto id needs
# to be shifted down by one to stay contiguous (cross references in columns
# I..M that point at labels rather than ids are left untouched, matching the
# regenerated sheet).
$lastRow = $ws.Range("A1048576").End(-4162).Row  # xlUp

for ($r = 73; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Value2
    if ($val -ne $null -and $val -match '^incentive:(\d+)$') {
        $num = [int]$Matches[1] - 1
        $cell.Value = "incentive:" + $num
    }
}
